$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1492.2
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H105").Value = 66931.664
$ws.Range("J105").Value = 66931.664
$ws.Range("L105").Value = 66931.664
$ws.Range("N105").Value = -73919.664
$ws.Range("H116").Value = 4250
$ws.Range("J116").Value = 4250
$ws.Range("L116").Value = 4250
$ws.Range("N116").Value = -11134
$ws.Range("H125").Value = 4547.364
$ws.Range("I125").Value = 1903.5
$ws.Range("K125").Value = 17131.5
$ws.Range("M125").Value = -14671.5
$ws.Range("H135").Value = 4032
$ws.Range("I135").Value = 4032
$ws.Range("K135").Value = 36288
$ws.Range("M135").Value = -33753

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5971.143
$ws.Range("I74").Value = 1574.5
$ws.Range("J74").Value = 11833.333
$ws.Range("K74").Value = 1574.5
$ws.Range("L74").Value = 11833.333
$ws.Range("M74").Value = -700.5
$ws.Range("N74").Value = -13581.333
$ws.Range("H77").Value = 5971.143
$ws.Range("I77").Value = 1574.5
$ws.Range("J77").Value = 11833.333
$ws.Range("K77").Value = 7872.5
$ws.Range("L77").Value = 59166.665
$ws.Range("M77").Value = -3504.5
$ws.Range("N77").Value = -67902.66500000001
$ws.Range("H94").Value = 24165
$ws.Range("J94").Value = 24165
$ws.Range("L94").Value = 24165
$ws.Range("N94").Value = -25967
$ws.Range("H98").Value = 46442.875
$ws.Range("J98").Value = 46442.875
$ws.Range("L98").Value = 46442.875
$ws.Range("N98").Value = -52432.875
$ws.Range("H130").Value = 19000
$ws.Range("J130").Value = 19000
$ws.Range("L130").Value = 19000
$ws.Range("N130").Value = -29040
$ws.Range("H132").Value = 4099.8
$ws.Range("I132").Value = 1749.5
$ws.Range("J132").Value = 5666.6665
$ws.Range("K132").Value = 5248.5
$ws.Range("L132").Value = 16999.9995
$ws.Range("M132").Value = -2718.5
$ws.Range("N132").Value = -22059.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 45
$ws.Range("I8").Value = 40
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 40
$ws.Range("L8").Value = 50
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = -330
$ws.Range("H11").Value = 750
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 750
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 750
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1030
$ws.Range("H103").Value = 5661.5
$ws.Range("J103").Value = 5661.5
$ws.Range("L103").Value = 5661.5
$ws.Range("N103").Value = -8005.5
$ws.Range("H105").Value = 6299.2
$ws.Range("I105").Value = 4874
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 4874
$ws.Range("L105").Value = 12000
$ws.Range("M105").Value = -3127
$ws.Range("N105").Value = -15494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 499.5
$ws.Range("I4").Value = 499.5
$ws.Range("K4").Value = 499.5
$ws.Range("M4").Value = -387.5
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H58").Value = 7324.3335
$ws.Range("I58").Value = 5486.625
$ws.Range("K58").Value = 5486.625
$ws.Range("M58").Value = -5283.625
$ws.Range("H99").Value = 10660
$ws.Range("J99").Value = 9500
$ws.Range("L99").Value = 9500
$ws.Range("N99").Value = -12496
$ws.Range("H118").Value = 29999
$ws.Range("J118").Value = 29999
$ws.Range("L118").Value = 29999
$ws.Range("N118").Value = -33313
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H126").Value = 10660
$ws.Range("J126").Value = 9500
$ws.Range("L126").Value = 28500
$ws.Range("N126").Value = -33440
$ws.Range("H136").Value = 7324.3335
$ws.Range("I136").Value = 5486.625
$ws.Range("K136").Value = 16459.875
$ws.Range("M136").Value = -13909.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200000720
$ws.Range("I4").Value = 897.5
$ws.Range("K4").Value = 2692.5
$ws.Range("M4").Value = -2580.5
$ws.Range("H5").Value = 1557
$ws.Range("I5").Value = 1483.1666
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 4449.4998
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -4337.4998
$ws.Range("N5").Value = -6224
$ws.Range("H11").Value = 1806.4
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 2253
$ws.Range("K11").Value = 60
$ws.Range("L11").Value = 6759
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = -7039
$ws.Range("H23").Value = 186.72728
$ws.Range("J23").Value = 126
$ws.Range("L23").Value = 378
$ws.Range("N23").Value = -848
$ws.Range("H34").Value = 816.6667
$ws.Range("J34").Value = 900
$ws.Range("L34").Value = 2700
$ws.Range("N34").Value = -2868
$ws.Range("H52").Value = 2021
$ws.Range("J52").Value = 2021
$ws.Range("L52").Value = 6063
$ws.Range("N52").Value = -6595
$ws.Range("H109").Value = 2454.3
$ws.Range("I109").Value = 2560.3333
$ws.Range("K109").Value = 7680.999899999999
$ws.Range("M109").Value = -6640.999899999999
$ws.Range("H135").Value = 1557
$ws.Range("I135").Value = 1483.1666
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 13348.4994
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -10813.4994
$ws.Range("N135").Value = -23070
$ws.Range("H139").Value = 1999.3636
$ws.Range("I139").Value = 1999.3
$ws.Range("J139").Value = 2000
$ws.Range("K139").Value = 5997.9
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = -857.8999999999996
$ws.Range("N139").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H99").Value = 11099.6
$ws.Range("I99").Value = 8874.75
$ws.Range("K99").Value = 8874.75
$ws.Range("M99").Value = -6628.75
$ws.Range("H102").Value = 3440.6365
$ws.Range("I102").Value = 3440.6365
$ws.Range("K102").Value = 3440.6365
$ws.Range("M102").Value = -1818.6365
$ws.Range("H132").Value = 7699.5
$ws.Range("J132").Value = 7319.6
$ws.Range("L132").Value = 21958.8
$ws.Range("N132").Value = -27018.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 764.8570999999999
$ws.Range("I2").Value = 470.8
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 470.8
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -358.8
$ws.Range("N2").Value = -1724
$ws.Range("H22").Value = 3083.3333
$ws.Range("I22").Value = 2700
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2700
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -2405
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 3083.3333
$ws.Range("I27").Value = 2700
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2700
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -2593
$ws.Range("N27").Value = -5214
$ws.Range("H122").Value = 3900
$ws.Range("J122").Value = 3900
$ws.Range("L122").Value = 11700
$ws.Range("N122").Value = -16600
$ws.Range("H132").Value = 7551.55
$ws.Range("J132").Value = 9174.5
$ws.Range("L132").Value = 27523.5
$ws.Range("N132").Value = -32583.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H132").Value = 4183.0835
$ws.Range("I132").Value = 2523.2354
$ws.Range("K132").Value = 7569.706200000001
$ws.Range("M132").Value = -5039.706200000001
